# Applies the scheduled-runner profit/price updates captured in the commit diff.
# Workbook has 8 sheets (one per crafting class); each hunk below updates the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for a single leve row.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # ALC
$ws2 = $wb.Worksheets.Item(2)  # ARM
$ws3 = $wb.Worksheets.Item(3)  # BSM
$ws4 = $wb.Worksheets.Item(4)  # CRP
$ws5 = $wb.Worksheets.Item(5)  # CUL
$ws6 = $wb.Worksheets.Item(6)  # GSM
$ws7 = $wb.Worksheets.Item(7)  # LTW
$ws8 = $wb.Worksheets.Item(8)  # WVR

# ALC row 31
$ws1.Range("H31").Value = 500
$ws1.Range("I31").Value = 500
$ws1.Range("K31").Value = 1500
$ws1.Range("M31").Value = -1270

# ALC row 33
$ws1.Range("H33").Value = 6062154
$ws1.Range("I33").Value = 1755.3077
$ws1.Range("K33").Value = 1755.3077
$ws1.Range("M33").Value = -1526.3077

# ALC row 58
$ws1.Range("H58").Value = 3889.5715
$ws1.Range("I58").Value = 2000
$ws1.Range("J58").Value = 4204.5
$ws1.Range("K58").Value = 6000
$ws1.Range("L58").Value = 12613.5
$ws1.Range("M58").Value = -5850
$ws1.Range("N58").Value = -12913.5

# ALC row 137
$ws1.Range("H137").Value = 1739.5366
$ws1.Range("I137").Value = 1251.1875
$ws1.Range("J137").Value = 2052.08
$ws1.Range("K137").Value = 3753.5625
$ws1.Range("L137").Value = 6156.24
$ws1.Range("M137").Value = -1203.5625
$ws1.Range("N137").Value = -11256.24

# ARM row 32
$ws2.Range("H32").Value = 4324.6235
$ws2.Range("I32").Value = 3730.195
$ws2.Range("K32").Value = 3730.195
$ws2.Range("M32").Value = -3443.195

# ARM row 92
$ws2.Range("H92").Value = 42000
$ws2.Range("J92").Value = 42000
$ws2.Range("L92").Value = 42000
$ws2.Range("N92").Value = -46992

# ARM row 132
$ws2.Range("H132").Value = 3295.6553
$ws2.Range("I132").Value = 1690.7727
$ws2.Range("J132").Value = 8339.571
$ws2.Range("K132").Value = 5072.3181
$ws2.Range("L132").Value = 25018.713
$ws2.Range("M132").Value = -2542.3181
$ws2.Range("N132").Value = -30078.713

# BSM row 20
$ws3.Range("H20").Value = 7960.3784
$ws3.Range("I20").Value = 1006.3077
$ws3.Range("K20").Value = 1006.3077
$ws3.Range("M20").Value = -759.3077

# BSM row 94
$ws3.Range("H94").Value = 2092.6
$ws3.Range("I94").Value = 1712.7142
$ws3.Range("J94").Value = 2425
$ws3.Range("K94").Value = 1712.7142
$ws3.Range("L94").Value = 2425
$ws3.Range("M94").Value = -1261.7142
$ws3.Range("N94").Value = -3327

# BSM row 107
$ws3.Range("H107").Value = 1024.4706
$ws3.Range("I107").Value = 956.0909
$ws3.Range("J107").Value = 1149.8334
$ws3.Range("K107").Value = 956.0909
$ws3.Range("L107").Value = 1149.8334
$ws3.Range("M107").Value = 963.9091
$ws3.Range("N107").Value = -4989.8334

# CRP row 2
$ws4.Range("H2").Value = 42401.8
$ws4.Range("I2").Value = 47001.332
$ws4.Range("J2").Value = 35502.5
$ws4.Range("K2").Value = 47001.332
$ws4.Range("L2").Value = 35502.5
$ws4.Range("M2").Value = -46888.332
$ws4.Range("N2").Value = -35728.5

# CRP row 11
$ws4.Range("H11").Value = 505
$ws4.Range("I11").Value = 505
$ws4.Range("K11").Value = 505
$ws4.Range("M11").Value = -365

# CRP row 58
$ws4.Range("H58").Value = 1466.15
$ws4.Range("I58").Value = 1162.1875
$ws4.Range("J58").Value = 2682
$ws4.Range("K58").Value = 1162.1875
$ws4.Range("L58").Value = 2682
$ws4.Range("M58").Value = -959.1875
$ws4.Range("N58").Value = -3088

# CRP row 132
$ws4.Range("H132").Value = 2588
$ws4.Range("I132").Value = 1467.1578
$ws4.Range("J132").Value = 5630.2856
$ws4.Range("K132").Value = 4401.4734
$ws4.Range("L132").Value = 16890.8568
$ws4.Range("M132").Value = -1871.4734
$ws4.Range("N132").Value = -21950.8568

# CRP row 134
$ws4.Range("H134").Value = 2634.1562
$ws4.Range("I134").Value = 2743.1
$ws4.Range("J134").Value = 1000
$ws4.Range("K134").Value = 8229.299999999999
$ws4.Range("L134").Value = 3000
$ws4.Range("M134").Value = -5694.299999999999
$ws4.Range("N134").Value = -8070

# CRP row 136
$ws4.Range("H136").Value = 1466.15
$ws4.Range("I136").Value = 1162.1875
$ws4.Range("J136").Value = 2682
$ws4.Range("K136").Value = 3486.5625
$ws4.Range("L136").Value = 8046
$ws4.Range("M136").Value = -936.5625
$ws4.Range("N136").Value = -13146

# CUL row 2
$ws5.Range("H2").Value = 1313.625
$ws5.Range("J2").Value = 95
$ws5.Range("L2").Value = 570
$ws5.Range("N2").Value = -796

# CUL row 17
$ws5.Range("H17").Value = 600
$ws5.Range("I17").Value = 600
$ws5.Range("J17").Value = 0
$ws5.Range("K17").Value = 1800
$ws5.Range("L17").Value = 0
$ws5.Range("M17").Value = -1631
$ws5.Range("N17").ClearContents()

# CUL row 113
$ws5.Range("H113").Value = 1200499.5
$ws5.Range("I113").Value = 1667126
$ws5.Range("K113").Value = 5001378
$ws5.Range("M113").Value = -4999208

# CUL row 129
$ws5.Range("H129").Value = 1964.6923
$ws5.Range("I129").Value = 1070
$ws5.Range("J129").Value = 2731.5715
$ws5.Range("K129").Value = 3210
$ws5.Range("L129").Value = 8194.7145
$ws5.Range("M129").Value = 1790
$ws5.Range("N129").Value = -18194.7145

# CUL row 131
$ws5.Range("H131").Value = 14286979
$ws5.Range("J131").Value = 14085824
$ws5.Range("L131").Value = 42257472
$ws5.Range("N131").Value = -42267552

# GSM row 80
$ws6.Range("H80").Value = 3862.5
$ws6.Range("I80").Value = 3780
$ws6.Range("K80").Value = 3780
$ws6.Range("M80").Value = -2782

# GSM row 83
$ws6.Range("H83").Value = 3862.5
$ws6.Range("I83").Value = 3780
$ws6.Range("K83").Value = 18900
$ws6.Range("M83").Value = -13908

# GSM row 97
$ws6.Range("H97").Value = 764.875
$ws6.Range("I97").Value = 753.1667
$ws6.Range("K97").Value = 753.1667
$ws6.Range("M97").Value = -257.1667

# LTW row 68
$ws7.Range("H68").Value = 142860660
$ws7.Range("J68").Value = 166670260
$ws7.Range("L68").Value = 166670260
$ws7.Range("N68").Value = -166671758

# LTW row 71
$ws7.Range("H71").Value = 142860660
$ws7.Range("J71").Value = 166670260
$ws7.Range("L71").Value = 833351300
$ws7.Range("N71").Value = -833358788

# WVR row 62
$ws8.Range("H62").Value = 8221.643
$ws8.Range("I62").Value = 4033.3333
$ws8.Range("J62").Value = 11362.875
$ws8.Range("K62").Value = 4033.3333
$ws8.Range("L62").Value = 11362.875
$ws8.Range("M62").Value = -3409.3333
$ws8.Range("N62").Value = -12610.875

# WVR row 65
$ws8.Range("H65").Value = 8221.643
$ws8.Range("I65").Value = 4033.3333
$ws8.Range("J65").Value = 11362.875
$ws8.Range("K65").Value = 20166.6665
$ws8.Range("L65").Value = 56814.375
$ws8.Range("M65").Value = -17046.6665
$ws8.Range("N65").Value = -63054.375

# WVR row 81
$ws8.Range("H81").Value = 2013
$ws8.Range("I81").Value = 2000
$ws8.Range("J81").Value = 2026
$ws8.Range("K81").Value = 4000
$ws8.Range("L81").Value = 4052
$ws8.Range("M81").Value = -2939
$ws8.Range("N81").Value = -6174

# WVR row 84
$ws8.Range("H84").Value = 2013
$ws8.Range("I84").Value = 2000
$ws8.Range("J84").Value = 2026
$ws8.Range("K84").Value = 20000
$ws8.Range("L84").Value = 20260
$ws8.Range("M84").Value = -14696
$ws8.Range("N84").Value = -30868

# WVR row 113
$ws8.Range("H113").Value = 872.7879
$ws8.Range("I113").Value = 676.1539
$ws8.Range("J113").Value = 1603.1428
$ws8.Range("K113").Value = 2028.4617
$ws8.Range("L113").Value = 4809.428400000001
$ws8.Range("M113").Value = 141.5382999999999
$ws8.Range("N113").Value = -9149.428400000001

# WVR row 132
$ws8.Range("H132").Value = 1474.2031
$ws8.Range("I132").Value = 946.73334
$ws8.Range("J132").Value = 2723.4736
$ws8.Range("K132").Value = 2840.20002
$ws8.Range("L132").Value = 8170.4208
$ws8.Range("M132").Value = -310.2000200000002
$ws8.Range("N132").Value = -13230.4208
